$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on D-column cells that will hold numeric-looking strings,
# so COM does not silently convert them to numbers, then restore default style.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "26.622.41"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.82%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.854.59"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.86%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "263.78"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.52%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.5258"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +1.28%  "

$ws.Range("E8").Value = "  +0.76%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06789"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("E10").Value = "  +0.67%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.7832"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +2.03%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.07778"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +1.22%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.848.92"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.54%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "88.53"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.54%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "5.029"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("E16").Value = "  -0.06%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "14.00"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.64%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.15%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.000007961"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.42%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "26.648.17"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.77%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "4.634"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +2.35%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "9.471"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "6.011"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.87%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "143.61"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.55%  "

$ws.Range("E25").Value = "  -6.67%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.696"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.94%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "17.05"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.96%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "111.86"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.81%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "4.183"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.40%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "4.108"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.27%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.08717"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.35%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.04862"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.63%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.132"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.22%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.7199"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +5.63%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.874"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.57%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "3.116"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +0.49%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.255"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +2.09%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01793"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +0.72%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.4871"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.78%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.9026"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +1.18%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "111.24"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("E42").Value = "  -3.20%  "

$ws.Range("E43").Value = "  +0.05%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "7.669"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("E45").Value = "  -0.06%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.05891"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.19%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "9.001"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.16%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "35.06"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.77%  "

$ws.Range("E49").Value = "  -1.84%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.8887"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +3.29%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "59.89"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.23%  "
